$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A92").NumberFormat = "@"
$ws.Range("A92").Value = "2025/10/11"
$ws.Range("A92").ClearFormats()

$ws.Range("B92").Value = "土"
$ws.Range("C92").Value = 6
$ws.Range("D92").Value = 201
